# Automatische test-sync: 2025-08-05 18:33:50
#
# - Logs sheet: append a new test-mail row (row 34)
# - Logs sheet: conditional formatting ranges grow from row 33 -> row 34
# - Dashboard sheet: swap category labels in rows 3/4
# - Dashboard sheet: append new "Documentatie / Datasheets" category (row 8)
# - Chart on Dashboard: category/value series ranges grow from row 7 -> row 8

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Logs sheet: add the new mail-log entry as row 34
# ---------------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A34").Value2 = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$wsLogs.Range("B34").Value2 = "mailmind.test@zohomail.eu"
$wsLogs.Range("C34").Value2 = "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?"
$wsLogs.Range("D34").Value2 = "Documentatie / Datasheets"
$wsLogs.Range("E34").Value2 = "Bedankt, we hebben dit doorgestuurd naar documentatie@bedrijf.nl."
$wsLogs.Range("F34").Value2 = "2025-08-05 18:32:50"
$wsLogs.Range("G34").Value2 = "Ja"
$wsLogs.Range("H34").Value2 = "Ja"
$wsLogs.Range("I34").Value2 = "Nee"
$wsLogs.Range("J34").Value2 = "Nee"

# Conditional formatting ranges need to grow by one row to cover the new entry
$wsLogs.Range("D2:D33").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D34"))
$wsLogs.Range("G2:G33").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G34"))
$wsLogs.Range("H2:H33").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("H2:H34"))
$wsLogs.Range("I2:I33").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("I2:I34"))
$wsLogs.Range("J2:J33").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("J2:J34"))

# ---------------------------------------------------------------------
# Dashboard sheet: swap rows 3/4 categories and append the new category
# ---------------------------------------------------------------------
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Range("A3").Value2 = "Inkoop / Bestellingen"
$wsDash.Range("A4").Value2 = "Klantenservice / Contact"

$wsDash.Range("A8").Value2 = "Documentatie / Datasheets"
$wsDash.Range("B8").Value2 = 1

# ---------------------------------------------------------------------
# Chart on Dashboard: extend the category/value series ranges to include
# the newly added row 8
# ---------------------------------------------------------------------
$chart = $wsDash.ChartObjects().Item(1).Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$8,'Dashboard'!`$B`$2:`$B`$8,1)"
